$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04026176048104
$ws.Range("D2").Value = 1.047921907626942
$ws.Range("E2").Value = 1.048651322513672
$ws.Range("F2").Value = 1.059437578510863
$ws.Range("I2").Value = 1.038449617506755
$ws.Range("J2").Value = 1.045349320086019
$ws.Range("K2").Value = 1.050683193543871
$ws.Range("L2").Value = 1.051410572982452
$ws.Range("M2").Value = 1.062167084062324
$ws.Range("N2").Value = 1.046833836023784

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041285047407972
$ws.Range("D3").Value = 1.048733097010849
$ws.Range("E3").Value = 1.049555622554827
$ws.Range("F3").Value = 1.060416330813621
$ws.Range("I3").Value = 1.038658904326011
$ws.Range("J3").Value = 1.046017783524972
$ws.Range("K3").Value = 1.05130632642991
$ws.Range("L3").Value = 1.052126722677755
$ws.Range("M3").Value = 1.062959649437724
$ws.Range("N3").Value = 1.047503248757494

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041947376886787
$ws.Range("D4").Value = 1.049257757592658
$ws.Range("E4").Value = 1.050141264632669
$ws.Range("F4").Value = 1.061050050049085
$ws.Range("I4").Value = 1.038792487284543
$ws.Range("J4").Value = 1.046449927868126
$ws.Range("K4").Value = 1.051708650735849
$ws.Range("L4").Value = 1.052589985795034
$ws.Range("M4").Value = 1.063472273764617
$ws.Range("N4").Value = 1.047936006795183

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042225866349044
$ws.Range("D5").Value = 1.049478268043926
$ws.Range("E5").Value = 1.050387587219633
$ws.Range("F5").Value = 1.061316560670008
$ws.Range("I5").Value = 1.038848204603456
$ws.Range("J5").Value = 1.046631505851609
$ws.Range("K5").Value = 1.051877575355392
$ws.Range("L5").Value = 1.052784708886888
$ws.Range("M5").Value = 1.063687727878705
$ws.Range("N5").Value = 1.04811784264023

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042272628649968
$ws.Range("D6").Value = 1.049515289354765
$ws.Range("E6").Value = 1.050428952797187
$ws.Range("F6").Value = 1.06136131452209
$ws.Range("I6").Value = 1.038857533912443
$ws.Range("J6").Value = 1.046661987975074
$ws.Range("K6").Value = 1.051905926056581
$ws.Range("L6").Value = 1.05281740181059
$ws.Range("M6").Value = 1.063723900430764
$ws.Range("N6").Value = 1.048148368051807

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04195109789577
$ws.Range("D7").Value = 1.04926070428882
$ws.Range("E7").Value = 1.050144555542052
$ws.Range("F7").Value = 1.061053610804759
$ws.Range("I7").Value = 1.038793233515445
$ws.Range("J7").Value = 1.046452354498092
$ws.Range("K7").Value = 1.051710908750962
$ws.Range("L7").Value = 1.052592587824024
$ws.Range("M7").Value = 1.063475152882968
$ws.Range("N7").Value = 1.047938436871242

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040607544598794
$ws.Range("D8").Value = 1.048196100543658
$ws.Range("E8").Value = 1.048956830973721
$ws.Range("F8").Value = 1.059768268820231
$ws.Range("I8").Value = 1.038520727540477
$ws.Range("J8").Value = 1.045575312083069
$ws.Range("K8").Value = 1.050893966845638
$ws.Range("L8").Value = 1.051652626353056
$ws.Range("M8").Value = 1.062434979966822
$ws.Range("N8").Value = 1.047060148955392

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038241537855794
$ws.Range("D9").Value = 1.046318391718991
$ws.Range("E9").Value = 1.046867774159797
$ws.Range("F9").Value = 1.057506445354947
$ws.Range("I9").Value = 1.038026471956487
$ws.Range("J9").Value = 1.044026845018561
$ws.Range("K9").Value = 1.049447672939156
$ws.Range("L9").Value = 1.049995302711203
$ws.Range("M9").Value = 1.060600426288212
$ws.Range("N9").Value = 1.04550948289004

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036665232440359
$ws.Range("D10").Value = 1.045065473963003
$ws.Range("E10").Value = 1.045477720963911
$ws.Range("F10").Value = 1.056000710466358
$ws.Range("I10").Value = 1.037687542740762
$ws.Range("J10").Value = 1.042992548559003
$ws.Range("K10").Value = 1.048478992308626
$ws.Range("L10").Value = 1.048889795051586
$ws.Range("M10").Value = 1.059376338554772
$ws.Range("N10").Value = 1.044473717610874

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035982922988514
$ws.Range("D11").Value = 1.044522694777824
$ws.Range("E11").Value = 1.044876451889901
$ws.Range("F11").Value = 1.055349230723346
$ws.Range("I11").Value = 1.037538553249134
$ws.Range("J11").Value = 1.042544222988155
$ws.Range("K11").Value = 1.048058488197563
$ws.Range("L11").Value = 1.048410957898306
$ws.Range("M11").Value = 1.058846055058167
$ws.Range("N11").Value = 1.044024755366292

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035729519277031
$ws.Range("D12").Value = 1.044321044605871
$ws.Range("E12").Value = 1.04465320968304
$ws.Range("F12").Value = 1.055107320063259
$ws.Range("I12").Value = 1.037482877189278
$ws.Range("J12").Value = 1.042377625037822
$ws.Range("K12").Value = 1.047902135794372
$ws.Range("L12").Value = 1.048233075260444
$ws.Range("M12").Value = 1.058649047975937
$ws.Range("N12").Value = 1.043857920827761

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035783873583625
$ws.Range("D13").Value = 1.044364300949504
$ws.Range("E13").Value = 1.044701091542157
$ws.Range("F13").Value = 1.055159207182168
$ws.Range("I13").Value = 1.03749483503329
$ws.Range("J13").Value = 1.04241336401759
$ws.Range("K13").Value = 1.04793568107557
$ws.Range("L13").Value = 1.048271232639927
$ws.Range("M13").Value = 1.058691308290011
$ws.Range("N13").Value = 1.043893710560979

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035961975815212
$ws.Range("D14").Value = 1.044506027078084
$ws.Range("E14").Value = 1.04485799665865
$ws.Range("F14").Value = 1.05532923271854
$ws.Range("I14").Value = 1.0375339578745
$ws.Range("J14").Value = 1.042530453372106
$ws.Range("K14").Value = 1.048045567287799
$ws.Range("L14").Value = 1.048396254493009
$ws.Range("M14").Value = 1.058829771125378
$ws.Range("N14").Value = 1.044010966195808

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036071715306242
$ws.Range("D15").Value = 1.04459334423864
$ws.Range("E15").Value = 1.044954683810158
$ws.Range("F15").Value = 1.055434001412275
$ws.Range("I15").Value = 1.037558018406431
$ws.Range("J15").Value = 1.042602586744005
$ws.Range("K15").Value = 1.048113250832607
$ws.Range("L15").Value = 1.048473281791524
$ws.Range("M15").Value = 1.058915077874443
$ws.Range("N15").Value = 1.04408320200537

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036710520167074
$ws.Range("D16").Value = 1.045101491078955
$ws.Range("E16").Value = 1.045517638617969
$ws.Range("F16").Value = 1.056043957915673
$ws.Range("I16").Value = 1.037697383717626
$ws.Range("J16").Value = 1.043022292624909
$ws.Range("K16").Value = 1.048506877525166
$ws.Range("L16").Value = 1.048921570891182
$ws.Range("M16").Value = 1.059411526612261
$ws.Range("N16").Value = 1.044503503916765

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037111290268624
$ws.Range("D17").Value = 1.045420169839157
$ws.Range("E17").Value = 1.045870935353771
$ws.Range("F17").Value = 1.05642670552711
$ws.Range("I17").Value = 1.037784206930288
$ws.Range("J17").Value = 1.043285437989224
$ws.Range("K17").Value = 1.048753506084814
$ws.Range("L17").Value = 1.049202732387081
$ws.Range("M17").Value = 1.059722870571124
$ws.Range("N17").Value = 1.044767022977681

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03734507594937
$ws.Range("D18").Value = 1.045606024862214
$ws.Range("E18").Value = 1.046077068398658
$ws.Range("F18").Value = 1.056650005216151
$ws.Range("I18").Value = 1.037834634162121
$ws.Range("J18").Value = 1.043438881003282
$ws.Range("K18").Value = 1.048897258129023
$ws.Range("L18").Value = 1.049366715096867
$ws.Range("M18").Value = 1.059904448508158
$ws.Range("N18").Value = 1.044920683898419

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037424794765476
$ws.Range("D19").Value = 1.045669392370908
$ws.Range("E19").Value = 1.046147364769087
$ws.Range("F19").Value = 1.056726153008553
$ws.Range("I19").Value = 1.037851792016511
$ws.Range("J19").Value = 1.043491193419661
$ws.Range("K19").Value = 1.048946256502727
$ws.Range("L19").Value = 1.049422626564759
$ws.Range("M19").Value = 1.05996635784671
$ws.Range("N19").Value = 1.04497307060443

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037068289022298
$ws.Range("D20").Value = 1.045385981193883
$ws.Range("E20").Value = 1.04583302363863
$ws.Range("F20").Value = 1.05638563520664
$ws.Range("I20").Value = 1.037774913885368
$ws.Range("J20").Value = 1.04325720965617
$ws.Range("K20").Value = 1.048727055742957
$ws.Range("L20").Value = 1.049172567866698
$ws.Range("M20").Value = 1.05968946875994
$ws.Range("N20").Value = 1.044738754557156

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035909528141006
$ws.Range("D21").Value = 1.044464293302536
$ws.Range("E21").Value = 1.044811789359098
$ws.Range("F21").Value = 1.055279162278448
$ws.Range("I21").Value = 1.037522446409522
$ws.Range("J21").Value = 1.042495975400399
$ws.Range("K21").Value = 1.048013212906031
$ws.Range("L21").Value = 1.04835943925382
$ws.Range("M21").Value = 1.058788998264683
$ws.Range("N21").Value = 1.043976439261427

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035181179856079
$ws.Range("D22").Value = 1.043884573184746
$ws.Range("E22").Value = 1.044170254057791
$ws.Range("F22").Value = 1.054583930063346
$ws.Range("I22").Value = 1.037361773663514
$ws.Range("J22").Value = 1.042016953198719
$ws.Range("K22").Value = 1.047563475137949
$ws.Range("L22").Value = 1.047848070662329
$ws.Range("M22").Value = 1.058222627734088
$ws.Range("N22").Value = 1.043496736793277

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035567271039428
$ws.Range("D23").Value = 1.044191914195047
$ws.Range("E23").Value = 1.04451029117185
$ws.Range("F23").Value = 1.054952442780321
$ws.Range("I23").Value = 1.037447132731135
$ws.Range("J23").Value = 1.042270930061766
$ws.Range("K23").Value = 1.047801976290678
$ws.Range("L23").Value = 1.048119168229591
$ws.Range("M23").Value = 1.058522891048264
$ws.Range("N23").Value = 1.0437510743326

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037087719362558
$ws.Range("D24").Value = 1.045401429648022
$ws.Range("E24").Value = 1.045850154122229
$ws.Range("F24").Value = 1.056404192963712
$ws.Range("I24").Value = 1.037779113677526
$ws.Range("J24").Value = 1.043269964964512
$ws.Range("K24").Value = 1.048739007829739
$ws.Range("L24").Value = 1.049186197958393
$ws.Range("M24").Value = 1.059704561673951
$ws.Range("N24").Value = 1.044751527979499

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038853026799866
$ws.Range("D25").Value = 1.046804024704104
$ws.Range("E25").Value = 1.047407381916874
$ws.Range("F25").Value = 1.058090806387619
$ws.Range("I25").Value = 1.038155912445243
$ws.Range("J25").Value = 1.044427513389067
$ws.Range("K25").Value = 1.049822368347069
$ws.Range("L25").Value = 1.050423873812586
$ws.Range("M25").Value = 1.061074891301995
$ws.Range("N25").Value = 1.045910720255589
